$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Versions")

$ws.Range("A18").Value = "17.0.0"
$ws.Range("B18").Value = "15.0.0"
$ws.Range("E18").Value = "5.1.0"

$ws.Range("A18").Font.Bold = $true

$ws.Range("A19").Select()
